$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.249.82"
$ws.Range("E2").Value = "  +0.75%  "

# Row 3
$ws.Range("D3").Value = "2.568.98"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "584.39"
$ws.Range("E5").Value = "  +3.13%  "

# Row 6
$ws.Range("D6").Value = "147.85"
$ws.Range("E6").Value = "  +1.34%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("E8").Value = "  +3.40%  "

# Row 9
$ws.Range("E9").Value = "  +3.92%  "

# Row 10
$ws.Range("E10").Value = "  +0.80%  "

# Row 11
$ws.Range("E11").Value = "  +0.41%  "

# Row 12
$ws.Range("E12").Value = "  +1.41%  "

# Row 13
$ws.Range("D13").Value = "27.49"
$ws.Range("E13").Value = "  +1.32%  "

# Row 14
$ws.Range("D14").Value = "3.028.72"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15
$ws.Range("D15").Value = "63.199.22"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16
$ws.Range("E16").Value = "  +4.30%  "

# Row 17
$ws.Range("D17").Value = "2.580.30"
$ws.Range("E17").Value = "  +1.45%  "

# Row 18
$ws.Range("D18").Value = "11.35"
$ws.Range("E18").Value = "  -0.78%  "

# Row 19
$ws.Range("D19").Value = "342.19"
$ws.Range("E19").Value = "  +2.55%  "

# Row 20
$ws.Range("E20").Value = "  +3.35%  "

# Row 21
$ws.Range("E21").Value = "  +2.00%  "

# Row 22
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").Value = "66.78"
$ws.Range("E23").Value = "  +3.34%  "

# Row 24
$ws.Range("D24").Value = "2.694.48"
$ws.Range("E24").Value = "  +1.59%  "

# Row 25
$ws.Range("E25").Value = "  +2.75%  "

# Row 26
$ws.Range("E26").Value = "  +1.22%  "

# Row 27
$ws.Range("D27").Value = "8.19"
$ws.Range("E27").Value = "  +12.82%  "

# Row 28
$ws.Range("D28").Value = "8.48"
$ws.Range("E28").Value = "  +2.00%  "

# Row 29 (was SuiNetwork, now Binance-PegBSC-USD)
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 (was Binance-PegBSC-USD, now SuiNetwork)
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "1.49"
$ws.Range("E30").Value = "  +0.20%  "

# Row 31
$ws.Range("E31").Value = "  +7.81%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0824"
$ws.Range("E32").Value = "  +2.28%  "

# Row 33
$ws.Range("D33").Value = "461.20"
$ws.Range("E33").Value = "  +13.30%  "

# Row 34
$ws.Range("D34").Value = "1.63"
$ws.Range("E34").Value = "  +3.78%  "

# Row 35
$ws.Range("D35").Value = "176.51"
$ws.Range("E35").Value = "  -0.12%  "

# Row 36
$ws.Range("E36").Value = "  +2.19%  "

# Row 37
$ws.Range("D37").Value = "19.23"
$ws.Range("E37").Value = "  +1.45%  "

# Row 38
$ws.Range("D38").Value = "4.52"
$ws.Range("E38").Value = "  +4.06%  "

# Row 39
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("E40").Value = "  -0.10%  "

# Row 41
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("D42").Value = "151.23"
$ws.Range("E42").Value = "  -1.05%  "

# Row 43
$ws.Range("D43").Value = "3.81"
$ws.Range("E43").Value = "  +2.13%  "

# Row 44
$ws.Range("D44").Value = "21.05"
$ws.Range("E44").Value = "  +2.01%  "

# Row 45
$ws.Range("D45").Value = "0.0552"
$ws.Range("E45").Value = "  +6.88%  "

# Row 46
$ws.Range("E46").Value = "  +1.83%  "

# Row 47
$ws.Range("E47").Value = "  +2.74%  "

# Row 48
$ws.Range("E48").Value = "  +2.18%  "

# Row 49
$ws.Range("D49").Value = "1.75"
$ws.Range("E49").Value = "  -0.65%  "

# Row 51
$ws.Range("E51").Value = "  +4.13%  "
